# Updated symbol list with refreshed Price (D) and Volume(1h) (E) values
# for the cryptos sheet, matching the GitHub Actions scrape commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of (cell, new value) pairs. Values are written as literal
# text, since the sheet stores Price/Volume as strings (e.g. "330.01",
# "0.03%"). Each cell is forced to Text format before the write so Excel
# does not auto-convert it to a number/percentage, then the style is
# restored to the default "Normal" style these cells originally had.
$updates = @(
    ,@("D2", "330.01")
    ,@("E2", "0.03%")
    ,@("D3", "43.82")
    ,@("E3", "6.78%")
    ,@("D4", "5.843")
    ,@("E4", "3.38%")
    ,@("D5", "0.08301")
    ,@("E5", "1.72%")
    ,@("D6", "8.798")
    ,@("E6", "0.56%")
    ,@("D7", "4.498")
    ,@("E7", "-1.02%")
    ,@("D8", "1.959")
    ,@("E8", "-3.76%")
    ,@("E9", "-1.59%")
    ,@("D10", "0.9296")
    ,@("E10", "1.21%")
    ,@("D11", "0.1251")
    ,@("E11", "-0.32%")
    ,@("D12", "0.1942")
    ,@("E12", "-0.55%")
    ,@("D13", "0.09507")
    ,@("E13", "1.91%")
    ,@("D14", "0.03973")
    ,@("E14", "7.59%")
    ,@("D15", "0.1063")
    ,@("E15", "0.75%")
    ,@("E16", "0.17%")
    ,@("D17", "0.005961")
    ,@("E17", "-3.41%")
    ,@("D18", "3.522")
    ,@("E18", "2.68%")
    ,@("D20", "9.038")
    ,@("E20", "9.25%")
    ,@("D21", "0.1370")
    ,@("E21", "-1.70%")
    ,@("D22", "0.2571")
    ,@("E22", "-3.13%")
    ,@("D23", "0.04393")
    ,@("E23", "-0.38%")
    ,@("D24", "0.001256")
    ,@("E24", "-1.11%")
    ,@("D25", "0.004406")
    ,@("E25", "2.46%")
    ,@("D26", "0.0001191")
    ,@("E26", "0.71%")
    ,@("D27", "0.0003991")
    ,@("E27", "-0.05%")
    ,@("D39", "0.02801")
    ,@("E39", "1.49%")
    ,@("E40", "2.75%")
    ,@("D41", "0.007910")
    ,@("E41", "3.30%")
    ,@("D42", "0.1423")
    ,@("E42", "0.56%")
    ,@("D43", "0.009075")
    ,@("E43", "-4.25%")
    ,@("D44", "0.002102")
    ,@("E44", "-0.61%")
    ,@("E45", "-13.06%")
    ,@("D46", "0.00007195")
    ,@("E46", "4.55%")
    ,@("D47", "0.00000000750")
    ,@("E47", "-0.16%")
    ,@("D48", "0.003792")
    ,@("E48", "5.90%")
    ,@("D49", "0.002279")
    ,@("E49", "-0.22%")
    ,@("D50", "0.00002101")
    ,@("E50", "-0.16%")
    ,@("D51", "0.0002001")
    ,@("E51", "-0.16%")
)

foreach ($pair in $updates) {
    $cell = $pair[0]
    $value = $pair[1]
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}
